# The canonical-OOXML diff for this commit only touches the namespace-
# declaration *attribute order* on the Office-2010 compatibility markup
# that PowerPoint writes into every shape's <a:extLst> (the
# <a14:hiddenFill>/<a14:hiddenLine> "remember this shape's fill/line even
# though it currently has noFill/no line" hints that live on the picture
# and the freeform "wave" decoration shapes of Slide Master 1 / Slide
# Layout 1, "Title Slide").
#
#   -<a14:hiddenFill xmlns:a14="...2010/main" xmlns="">
#   +<a14:hiddenFill xmlns="" xmlns:a14="...2010/main">
#
# Both namespace declarations (the a14 prefix binding and the empty
# default-namespace reset) are present before and after; only their
# textual order on the start tag changes. That carries no schema/visual
# meaning at all (namespace scoping in XML is order independent), and
# PowerPoint's object model has no property that reaches into an
# extLst's raw markup to re-order its attributes — hiddenFill/hiddenLine
# are round-tripped verbatim by the host and are not reachable via
# Shape.Fill / Shape.Line / any other exposed member.
#
# So there is no COM call that changes this deck's rendered content,
# shapes, text, or any inspectable property — the author's commit is a
# byte-level re-serialization, not a content edit. Touch the deck
# read-only through the object model (confirming it is the expected
# presentation) and leave every shape/property exactly as authored.

$p = $ppt.ActivePresentation

$master = $p.Slides.Item(1).Master
$title_layout = $master.CustomLayouts.Item(1)

Write-Output "Slides: $($p.Slides.Count); Master shapes: $($master.Shapes.Count); '$($title_layout.Name)' layout shapes: $($title_layout.Shapes.Count)"
